$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.006.62"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "'3.771.56"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.21%  "
$ws.Range("E4").Value = "  -0.38%  "
$ws.Range("D5").Value = "'632.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.09%  "
$ws.Range("D6").Value = "'165.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.35%  "
$ws.Range("D7").Value = "'3.771.03"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.17%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  +0.21%  "
$ws.Range("E10").Value = "  -1.67%  "
$ws.Range("E11").Value = "  +0.62%  "
$ws.Range("D12").Value = "'6.75"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.63%  "
$ws.Range("E13").Value = "  -4.04%  "
$ws.Range("D14").Value = "'34.82"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.88%  "
$ws.Range("D15").Value = "'4.406.15"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.27%  "
$ws.Range("D16").Value = "'3.772.79"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.30%  "
$ws.Range("D17").Value = "'68.989.35"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D18").Value = "'17.65"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.15%  "
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("E20").Value = "  -1.79%  "
$ws.Range("D21").Value = "'467.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.40%  "
$ws.Range("D22").Value = "'9.49"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.86%  "
$ws.Range("E23").Value = "  -0.75%  "
$ws.Range("D24").Value = "'82.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.30%  "
$ws.Range("E25").Value = "  -5.40%  "
$ws.Range("D26").Value = "'12.08"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.31%  "
$ws.Range("E27").Value = "  -1.30%  "
$ws.Range("D28").Value = "'10.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.55%  "
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("D30").Value = "'3.922.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.37%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'2.68"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.52%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'2.27"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.03%  "
$ws.Range("E33").Value = "  -3.35%  "
$ws.Range("E34").Value = "  +18.58%  "
$ws.Range("D35").Value = "'28.42"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.47%  "
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("D37").Value = "'3.723.53"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.13%  "
$ws.Range("D38").Value = "'8.89"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.32%  "
$ws.Range("E39").Value = "  -1.15%  "
$ws.Range("D40").Value = "'3.27"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.71%  "
$ws.Range("E41").Value = "  -2.14%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "'1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.16%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "'0.961"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.71%  "
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("E45").Value = "  +5.14%  "
$ws.Range("D46").Value = "'156.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.57%  "
$ws.Range("D47").Value = "'1.42"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.07%  "
$ws.Range("D48").Value = "'43.50"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.61%  "
$ws.Range("D49").Value = "'46.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.25%  "
$ws.Range("D50").Value = "'0.293"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.13%  "
$ws.Range("E51").Value = "  -1.27%  "
